$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 17508.5
$ws.Range("I21").Value = 18017
$ws.Range("J21").Value = 17000
$ws.Range("K21").Value = 18017
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = -17549
$ws.Range("N21").Value = -17936

# Hunk 1: sheet ALC, row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 17508.5
$ws.Range("I23").Value = 18017
$ws.Range("J23").Value = 17000
$ws.Range("K23").Value = 18017
$ws.Range("L23").Value = 17000
$ws.Range("M23").Value = -17783
$ws.Range("N23").Value = -17468

# Hunk 2: sheet ALC, row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 587.5
$ws.Range("I29").Value = 125
$ws.Range("J29").Value = 1050
$ws.Range("K29").Value = 375
$ws.Range("L29").Value = 3150
$ws.Range("M29").Value = -94
$ws.Range("N29").Value = -3712

# Hunk 3: sheet ALC, row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1617.7241
$ws.Range("I38").Value = 238.84616
$ws.Range("J38").Value = 2738.0625
$ws.Range("K38").Value = 716.5384799999999
$ws.Range("L38").Value = 8214.1875
$ws.Range("M38").Value = -344.5384799999999
$ws.Range("N38").Value = -8958.1875

# Hunk 4: sheet ALC, row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 348
$ws.Range("I58").Value = 316.375
$ws.Range("J58").Value = 432.33334
$ws.Range("K58").Value = 949.125
$ws.Range("L58").Value = 1297.00002
$ws.Range("M58").Value = -799.125
$ws.Range("N58").Value = -1597.00002

# Hunk 5: sheet ALC, row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 15486.818
$ws.Range("J87").Value = 15486.818
$ws.Range("L87").Value = 15486.818
$ws.Range("N87").Value = -17982.818

# Hunk 6: sheet ALC, row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 15486.818
$ws.Range("J90").Value = 15486.818
$ws.Range("L90").Value = 46460.454
$ws.Range("N90").Value = -58940.454

# Hunk 7: sheet ARM, row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 30006
$ws.Range("J23").Value = 30006
$ws.Range("L23").Value = 30006
$ws.Range("N23").Value = -30524

# Hunk 8: sheet ARM, row 37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 27520.75
$ws.Range("I37").Value = 9900
$ws.Range("J37").Value = 30038
$ws.Range("K37").Value = 9900
$ws.Range("L37").Value = 30038
$ws.Range("M37").Value = -9627
$ws.Range("N37").Value = -30584

# Hunk 9: sheet ARM, row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 30024.166
$ws.Range("J44").Value = 30024.166
$ws.Range("L44").Value = 30024.166
$ws.Range("N44").Value = -31000.166

# Hunk 10: sheet ARM, row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 29500
$ws.Range("J55").Value = 29500
$ws.Range("L55").Value = 29500
$ws.Range("N55").Value = -30130

# Hunk 11: sheet ARM, row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3209.3635
$ws.Range("I63").Value = 3209.3635
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3209.3635
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2523.3635
$ws.Range("N63").ClearContents()

# Hunk 12: sheet ARM, row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3209.3635
$ws.Range("I66").Value = 3209.3635
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 16046.8175
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -12614.8175
$ws.Range("N66").ClearContents()

# Hunk 13: sheet ARM, row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 29504.545
$ws.Range("J80").Value = 29504.545
$ws.Range("L80").Value = 29504.545
$ws.Range("N80").Value = -31500.545

# Hunk 14: sheet ARM, row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 29504.545
$ws.Range("J83").Value = 29504.545
$ws.Range("L83").Value = 88513.63499999999
$ws.Range("N83").Value = -98497.63499999999

# Hunk 15: sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1426.3704
$ws.Range("I132").Value = 1407.909
$ws.Range("J132").Value = 1465.4231
$ws.Range("K132").Value = 4223.727000000001
$ws.Range("L132").Value = 4396.2693
$ws.Range("M132").Value = -1693.727000000001
$ws.Range("N132").Value = -9456.2693

# Hunk 16: sheet BSM, row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 11000
$ws.Range("J35").Value = 11000
$ws.Range("L35").Value = 11000
$ws.Range("N35").Value = -11620

# Hunk 17: sheet BSM, row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13909.091
$ws.Range("I82").Value = 2571.4285
$ws.Range("J82").Value = 33750
$ws.Range("K82").Value = 2571.4285
$ws.Range("L82").Value = 33750
$ws.Range("M82").Value = -2188.4285
$ws.Range("N82").Value = -34516

# Hunk 18: sheet BSM, row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 13909.091
$ws.Range("I85").Value = 2571.4285
$ws.Range("J85").Value = 33750
$ws.Range("K85").Value = 2571.4285
$ws.Range("L85").Value = 33750
$ws.Range("M85").Value = -1245.4285
$ws.Range("N85").Value = -36402

# Hunk 19: sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28194
$ws.Range("I31").Value = 34865.613
$ws.Range("J31").Value = 16028.117
$ws.Range("K31").Value = 34865.613
$ws.Range("L31").Value = 16028.117
$ws.Range("M31").Value = -34570.613
$ws.Range("N31").Value = -16618.117

# Hunk 20: sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 28194
$ws.Range("I34").Value = 34865.613
$ws.Range("J34").Value = 16028.117
$ws.Range("K34").Value = 34865.613
$ws.Range("L34").Value = 16028.117
$ws.Range("M34").Value = -34663.613
$ws.Range("N34").Value = -16432.117

# Hunk 21: sheet CRP, row 39
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 13670
$ws.Range("I39").Value = 16015
$ws.Range("J39").Value = 6635
$ws.Range("K39").Value = 16015
$ws.Range("L39").Value = 6635
$ws.Range("M39").Value = -15624
$ws.Range("N39").Value = -7417

# Hunk 22: sheet CRP, row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 14548
$ws.Range("I41").Value = 5353
$ws.Range("K41").Value = 5353
$ws.Range("M41").Value = -4925

# Hunk 23: sheet CRP, row 49
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 13670
$ws.Range("I49").Value = 16015
$ws.Range("J49").Value = 6635
$ws.Range("K49").Value = 16015
$ws.Range("L49").Value = 6635
$ws.Range("M49").Value = -15833
$ws.Range("N49").Value = -6999

# Hunk 24: sheet CRP, row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 21342
$ws.Range("J50").Value = 21342
$ws.Range("L50").Value = 21342
$ws.Range("N50").Value = -22592

# Hunk 25: sheet CRP, row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 33934.08
$ws.Range("J59").Value = 33934.08
$ws.Range("L59").Value = 33934.08
$ws.Range("N59").Value = -36224.08

# Hunk 26: sheet CRP, row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 17010
$ws.Range("I60").Value = 2596.5
$ws.Range("J60").Value = 19892.7
$ws.Range("K60").Value = 2596.5
$ws.Range("L60").Value = 19892.7
$ws.Range("M60").Value = -2085.5
$ws.Range("N60").Value = -20914.7

# Hunk 27: sheet CRP, row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 29921.166
$ws.Range("J74").Value = 29921.166
$ws.Range("L74").Value = 29921.166
$ws.Range("N74").Value = -31669.166

# Hunk 28: sheet CRP, row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 29921.166
$ws.Range("J77").Value = 29921.166
$ws.Range("L77").Value = 89763.49800000001
$ws.Range("N77").Value = -98499.49800000001

# Hunk 29: sheet CUL, row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1606.8823
$ws.Range("J34").Value = 1458.2
$ws.Range("L34").Value = 4374.6
$ws.Range("N34").Value = -4542.6

# Hunk 30: sheet CUL, row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2093.2173
$ws.Range("J39").Value = 2093.2173
$ws.Range("L39").Value = 6279.651899999999
$ws.Range("N39").Value = -6867.651899999999

# Hunk 31: sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2236.5293
$ws.Range("J55").Value = 2909.5454
$ws.Range("L55").Value = 8728.636200000001
$ws.Range("N55").Value = -9082.636200000001

# Hunk 32: sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 792.36536
$ws.Range("I131").Value = 432.66666
$ws.Range("J131").Value = 938.1892
$ws.Range("K131").Value = 1297.99998
$ws.Range("L131").Value = 2814.5676
$ws.Range("M131").Value = 3742.00002
$ws.Range("N131").Value = -12894.5676

# Hunk 33: sheet GSM, row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 16262.846
$ws.Range("I43").Value = 3805.6667
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 3805.6667
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = -3654.6667
$ws.Range("N43").Value = -20302

# Hunk 34: sheet GSM, row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 17753.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 17753.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 17753.75
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -18065.75

# Hunk 35: sheet GSM, row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16659.8
$ws.Range("I57").Value = 3055
$ws.Range("J57").Value = 20061
$ws.Range("K57").Value = 3055
$ws.Range("L57").Value = 20061
$ws.Range("M57").Value = -2235
$ws.Range("N57").Value = -21701

# Hunk 36: sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1984.6875
$ws.Range("I80").Value = 1867.0834
$ws.Range("J80").Value = 2337.5
$ws.Range("K80").Value = 1867.0834
$ws.Range("L80").Value = 2337.5
$ws.Range("M80").Value = -869.0834
$ws.Range("N80").Value = -4333.5

# Hunk 37: sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1984.6875
$ws.Range("I83").Value = 1867.0834
$ws.Range("J83").Value = 2337.5
$ws.Range("K83").Value = 9335.416999999999
$ws.Range("L83").Value = 11687.5
$ws.Range("M83").Value = -4343.416999999999
$ws.Range("N83").Value = -21671.5

# Hunk 38: sheet WVR, row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30661.6
$ws.Range("J54").Value = 30661.6
$ws.Range("L54").Value = 30661.6
$ws.Range("N54").Value = -31701.6

# Hunk 39: sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1286.375
$ws.Range("I81").Value = 1286.375
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2572.75
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1511.75
$ws.Range("N81").ClearContents()

# Hunk 40: sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1286.375
$ws.Range("I84").Value = 1286.375
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 12863.75
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -7559.75
$ws.Range("N84").ClearContents()

# Hunk 41: sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 352.69232
$ws.Range("I136").Value = 252.46153
$ws.Range("J136").Value = 653.38464
$ws.Range("K136").Value = 757.38459
$ws.Range("L136").Value = 1960.15392
$ws.Range("M136").Value = 1792.61541
$ws.Range("N136").Value = -7060.15392
